# Commit 08: Add Folder And Second Chap
#
# Marks the "2. Thiet ke he thong" (system design) chapter's tasks
# (rows 12-20, column F / TINH TRANG) as "Xong" (Done) - the same way the
# first chapter's tasks (rows 7-10) are already marked - and nudges the
# saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectPlan")

# Column F ("TINH TRANG") for every task row under chapter 2 gets set to
# "Xong", matching the formatting already used on F7:F10 (bordered cell,
# horizontally centered).
$doneRows = 12,13,14,15,16,17,18,19,20
foreach ($r in $doneRows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = "Xong"
    $cell.HorizontalAlignment = -4108
}

# Restore the saved scroll position / active selection for the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H17").Select()
